$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.917.29'
$ws.Range('E2').Value = '  +2.12%  '
$ws.Range('D3').Value = '3.808.98'
$ws.Range('E3').Value = '  +0.86%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '630.96'
$ws.Range('E5').Value = '  +5.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '164.81'
$ws.Range('E6').Value = '  +0.24%  '
$ws.Range('D7').Value = '3.806.82'
$ws.Range('E7').Value = '  +0.85%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.519'
$ws.Range('E9').Value = '  +0.93%  '
$ws.Range('E10').Value = '  +2.72%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.452'
$ws.Range('E11').Value = '  +0.73%  '
$ws.Range('E12').Value = '  +3.39%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000249'
$ws.Range('E13').Value = '  +0.65%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.89'
$ws.Range('E14').Value = '  +1.19%  '
$ws.Range('D15').Value = '4.449.66'
$ws.Range('E15').Value = '  +0.85%  '
$ws.Range('D16').Value = '3.782.76'
$ws.Range('E16').Value = '  -0.03%  '
$ws.Range('D17').Value = '68.983.62'
$ws.Range('E17').Value = '  +2.09%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.95'
$ws.Range('E18').Value = '  -1.60%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.11'
$ws.Range('E19').Value = '  +1.24%  '
$ws.Range('E20').Value = '  -0.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '465.46'
$ws.Range('E21').Value = '  +1.32%  '
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.707'
$ws.Range('E23').Value = '  +1.89%  '
$ws.Range('E24').Value = '  +4.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.54'
$ws.Range('E25').Value = '  +1.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.98'
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('E27').Value = '  +2.80%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.05'
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('D30').Value = '3.962.45'
$ws.Range('E30').Value = '  +0.92%  '
$ws.Range('E31').Value = '  +4.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.21'
$ws.Range('E32').Value = '  +1.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.26'
$ws.Range('E33').Value = '  -2.29%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '29.18'
$ws.Range('E34').Value = '  +0.57%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.04'
$ws.Range('E36').Value = '  +0.99%  '
$ws.Range('E37').Value = '  +3.71%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.149'
$ws.Range('E38').Value = '  +7.98%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.42'
$ws.Range('E39').Value = '  +5.70%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.89'
$ws.Range('E40').Value = '  +2.63%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.976'
$ws.Range('E41').Value = '  -0.96%  '
$ws.Range('E42').Value = '  +0.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '157.69'
$ws.Range('E44').Value = '  +4.04%  '
$ws.Range('E45').Value = '  +6.19%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.299'
$ws.Range('E46').Value = '  +0.95%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '46.84'
$ws.Range('E47').Value = '  -1.10%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '42.67'
$ws.Range('E48').Value = '  -1.77%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.42'
$ws.Range('E49').Value = '  +1.34%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.89'
$ws.Range('E50').Value = '  +2.99%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.000280'
$ws.Range('E51').Value = '  +14.08%  '
